$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 226, shifting the existing rows 226..342 down to 228..344
$ws.Rows.Item(226).Insert()
$ws.Rows.Item(226).Insert()

# Populate the first new row (226)
$ws.Range("A226").Value = 6
$ws.Range("B226").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C226").Value = "Metropolitana"
$ws.Range("D226").Value = 44452
$ws.Range("E226").Value = 13
$ws.Range("F226").Value = 100112012
$ws.Range("G226").Value = "Espinaca"
$ws.Range("H226").Value = "Sin especificar"
$ws.Range("I226").Value = "Primera"
$ws.Range("J226").Value = 380
$ws.Range("K226").Value = 4500
$ws.Range("L226").Value = 5000
$ws.Range("M226").Value = 4697
$ws.Range("N226").Value = "`$/cuna 10 kilos"
$ws.Range("O226").Value = "Provincia de Chacabuco"
$ws.Range("P226").Value = 470
$ws.Range("Q226").Value = 10
$ws.Range("R226").Value = "Hortaliza"

# Populate the second new row (227)
$ws.Range("A227").Value = 6
$ws.Range("B227").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C227").Value = "Metropolitana"
$ws.Range("D227").Value = 44452
$ws.Range("E227").Value = 13
$ws.Range("F227").Value = 100112012
$ws.Range("G227").Value = "Espinaca"
$ws.Range("H227").Value = "Sin especificar"
$ws.Range("I227").Value = "Primera"
$ws.Range("J227").Value = 450
$ws.Range("K227").Value = 4500
$ws.Range("L227").Value = 5000
$ws.Range("M227").Value = 4711
$ws.Range("N227").Value = "`$/cuna 10 kilos"
$ws.Range("O227").Value = "Región Metropolitana"
$ws.Range("P227").Value = 471
$ws.Range("Q227").Value = 10
$ws.Range("R227").Value = "Hortaliza"
